$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3
$ws.Range("C2").Value = 0.3
$ws.Range("D2").Value = 0.3
$ws.Range("B3").Value = 0.5
$ws.Range("C3").Value = 0.5
$ws.Range("D3").Value = 0.5
$ws.Range("B4").Value = 0.4166666666666667
$ws.Range("C4").Value = 0.4166666666666667
$ws.Range("D4").Value = 0.4166666666666667
$ws.Range("E4").Value = 0.4166666666666667
$ws.Range("B5").Value = 0.4
$ws.Range("C5").Value = 0.4
$ws.Range("D5").Value = 0.4
$ws.Range("B6").Value = 0.4166666666666667
$ws.Range("C6").Value = 0.4166666666666667
$ws.Range("D6").Value = 0.4166666666666667
$ws.Range("B7").Value = 0.375
$ws.Range("C7").Value = 0.6
$ws.Range("D7").Value = 0.4615384615384615
$ws.Range("B8").Value = 0.5
$ws.Range("C8").Value = 0.2857142857142857
$ws.Range("D8").Value = 0.3636363636363636
$ws.Range("B9").Value = 0.4166666666666667
$ws.Range("C9").Value = 0.4166666666666667
$ws.Range("D9").Value = 0.4166666666666667
$ws.Range("E9").Value = 0.4166666666666667
$ws.Range("B10").Value = 0.4375
$ws.Range("C10").Value = 0.4428571428571428
$ws.Range("D10").Value = 0.4125874125874126
$ws.Range("B11").Value = 0.4479166666666667
$ws.Range("C11").Value = 0.4166666666666667
$ws.Range("D11").Value = 0.4044289044289044
$ws.Range("B12").Value = 0.4285714285714285
$ws.Range("D12").Value = 0.5
$ws.Range("B13").Value = 0.6
$ws.Range("C13").Value = 0.4285714285714285
$ws.Range("D13").Value = 0.5
$ws.Range("B14").Value = 0.5
$ws.Range("C14").Value = 0.5
$ws.Range("D14").Value = 0.5
$ws.Range("E14").Value = 0.5
$ws.Range("B15").Value = 0.5142857142857142
$ws.Range("C15").Value = 0.5142857142857142
$ws.Range("D15").Value = 0.5
$ws.Range("B16").Value = 0.5285714285714286
$ws.Range("C16").Value = 0.5
$ws.Range("D16").Value = 0.5
$ws.Range("B17").Value = 0.5
$ws.Range("D17").Value = 0.5
$ws.Range("B18").Value = 0.6428571428571429
$ws.Range("C18").Value = 0.6428571428571429
$ws.Range("D18").Value = 0.6428571428571429
$ws.Range("B19").Value = 0.5833333333333334
$ws.Range("C19").Value = 0.5833333333333334
$ws.Range("D19").Value = 0.5833333333333334
$ws.Range("E19").Value = 0.5833333333333334
$ws.Range("B20").Value = 0.5714285714285714
$ws.Range("C20").Value = 0.5714285714285714
$ws.Range("D20").Value = 0.5714285714285714
$ws.Range("B21").Value = 0.5833333333333334
$ws.Range("C21").Value = 0.5833333333333334
$ws.Range("D21").Value = 0.5833333333333334
$ws.Range("B22").Value = 0.3846153846153846
$ws.Range("C22").Value = 0.5
$ws.Range("D22").Value = 0.4347826086956522
$ws.Range("B23").Value = 0.5454545454545454
$ws.Range("C23").Value = 0.4285714285714285
$ws.Range("D23").Value = 0.4799999999999999
$ws.Range("B24").Value = 0.4583333333333333
$ws.Range("C24").Value = 0.4583333333333333
$ws.Range("D24").Value = 0.4583333333333333
$ws.Range("E24").Value = 0.4583333333333333
$ws.Range("B25").Value = 0.465034965034965
$ws.Range("C25").Value = 0.4642857142857143
$ws.Range("D25").Value = 0.457391304347826
$ws.Range("B26").Value = 0.4784382284382284
$ws.Range("C26").Value = 0.4583333333333333
$ws.Range("D26").Value = 0.4611594202898551
